$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 128.42857
$ws.Range("I31").Value = 128.42857
$ws.Range("K31").Value = 385.28571
$ws.Range("M31").Value = -155.28571
$ws.Range("H49").Value = 3814
$ws.Range("J49").Value = 3728
$ws.Range("L49").Value = 11184
$ws.Range("N49").Value = -11456
$ws.Range("H55").Value = 160.90909
$ws.Range("I55").Value = 147
$ws.Range("K55").Value = 147
$ws.Range("M55").Value = 67
$ws.Range("H69").Value = 17184.375
$ws.Range("J69").Value = 17610.727
$ws.Range("L69").Value = 52832.181
$ws.Range("N69").Value = -54580.181
$ws.Range("H72").Value = 17184.375
$ws.Range("J72").Value = 17610.727
$ws.Range("L72").Value = 158496.543
$ws.Range("N72").Value = -167232.543
$ws.Range("H100").Value = 1710
$ws.Range("I100").Value = 1805
$ws.Range("J100").Value = 1266.6666
$ws.Range("K100").Value = 1805
$ws.Range("L100").Value = 1266.6666
$ws.Range("M100").Value = -1264
$ws.Range("N100").Value = -2348.6666
$ws.Range("H116").Value = 7572.4375
$ws.Range("I116").Value = 7074.8887
$ws.Range("K116").Value = 7074.8887
$ws.Range("M116").Value = -3632.8887
$ws.Range("H141").Value = 5465
$ws.Range("I141").Value = 4225.1665
$ws.Range("J141").Value = 7324.75
$ws.Range("K141").Value = 12675.4995
$ws.Range("L141").Value = 21974.25
$ws.Range("M141").Value = -7495.499500000002
$ws.Range("N141").Value = -32334.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6628.421
$ws.Range("I74").Value = 3394.8333
$ws.Range("K74").Value = 3394.8333
$ws.Range("M74").Value = -2520.8333
$ws.Range("H77").Value = 6628.421
$ws.Range("I77").Value = 3394.8333
$ws.Range("K77").Value = 16974.1665
$ws.Range("M77").Value = -12606.1665
$ws.Range("H97").Value = 350.47058
$ws.Range("I97").Value = 256.66666
$ws.Range("K97").Value = 256.66666
$ws.Range("M97").Value = 239.33334
$ws.Range("H102").Value = 1474.3077
$ws.Range("I102").Value = 1055.6364
$ws.Range("K102").Value = 1055.6364
$ws.Range("M102").Value = 566.3635999999999
$ws.Range("H110").Value = 2183.6667
$ws.Range("I110").Value = 2241.5293
$ws.Range("J110").Value = 1200
$ws.Range("K110").Value = 2241.5293
$ws.Range("L110").Value = 1200
$ws.Range("M110").Value = -196.5293000000001
$ws.Range("N110").Value = -5290

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 31466.5
$ws.Range("J21").Value = 31466.5
$ws.Range("L21").Value = 31466.5
$ws.Range("N21").Value = -31938.5
$ws.Range("H22").Value = 637.7692
$ws.Range("I22").Value = 578.26666
$ws.Range("J22").Value = 718.9091
$ws.Range("K22").Value = 578.26666
$ws.Range("L22").Value = 718.9091
$ws.Range("M22").Value = -405.26666
$ws.Range("N22").Value = -1064.9091
$ws.Range("H54").Value = 8020.25
$ws.Range("J54").Value = 8499
$ws.Range("L54").Value = 8499
$ws.Range("N54").Value = -9467
$ws.Range("H94").Value = 813.35
$ws.Range("I94").Value = 821.4211
$ws.Range("J94").Value = 660
$ws.Range("K94").Value = 821.4211
$ws.Range("L94").Value = 660
$ws.Range("M94").Value = -370.4211
$ws.Range("N94").Value = -1562
$ws.Range("H105").Value = 1305.0286
$ws.Range("I105").Value = 1314.3334
$ws.Range("K105").Value = 1314.3334
$ws.Range("M105").Value = 432.6666
$ws.Range("H134").Value = 3315.1667
$ws.Range("I134").Value = 1978.2
$ws.Range("K134").Value = 5934.6
$ws.Range("M134").Value = -3399.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1110.4762
$ws.Range("I16").Value = 915.64703
$ws.Range("J16").Value = 1938.5
$ws.Range("K16").Value = 915.64703
$ws.Range("L16").Value = 1938.5
$ws.Range("M16").Value = -628.64703
$ws.Range("N16").Value = -2512.5
$ws.Range("H38").Value = 16000
$ws.Range("I38").Value = 6000
$ws.Range("K38").Value = 6000
$ws.Range("M38").Value = -5623
$ws.Range("H39").Value = 23274.572
$ws.Range("I39").Value = 18204.4
$ws.Range("K39").Value = 18204.4
$ws.Range("M39").Value = -17813.4
$ws.Range("H46").Value = 16000
$ws.Range("I46").Value = 6000
$ws.Range("K46").Value = 6000
$ws.Range("M46").Value = -5789
$ws.Range("H48").Value = 50000
$ws.Range("J48").Value = 50000
$ws.Range("L48").Value = 50000
$ws.Range("N48").Value = -50952
$ws.Range("H49").Value = 23274.572
$ws.Range("I49").Value = 18204.4
$ws.Range("K49").Value = 18204.4
$ws.Range("M49").Value = -18022.4
$ws.Range("H86").Value = 7325
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 7325
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H113").Value = 1110.4762
$ws.Range("I113").Value = 915.64703
$ws.Range("J113").Value = 1938.5
$ws.Range("K113").Value = 915.64703
$ws.Range("L113").Value = 1938.5
$ws.Range("M113").Value = 1254.35297
$ws.Range("N113").Value = -6278.5
$ws.Range("H139").Value = 98419.75
$ws.Range("J139").Value = 98419.75
$ws.Range("L139").Value = 98419.75
$ws.Range("N139").Value = -108699.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 25001740
$ws.Range("I4").Value = 89286670
$ws.Range("K4").Value = 267860010
$ws.Range("M4").Value = -267859898
$ws.Range("H13").Value = 201
$ws.Range("I13").Value = 152.5
$ws.Range("K13").Value = 457.5
$ws.Range("M13").Value = -289.5
$ws.Range("H51").Value = 134
$ws.Range("I51").Value = 134
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 402
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = 58
$ws.Range("H113").Value = 1428.4445
$ws.Range("J113").Value = 2260.4
$ws.Range("L113").Value = 6781.200000000001
$ws.Range("N113").Value = -11121.2
$ws.Range("H129").Value = 18543166
$ws.Range("I129").Value = 18636
$ws.Range("J129").Value = 55592224
$ws.Range("K129").Value = 55908
$ws.Range("L129").Value = 166776672
$ws.Range("M129").Value = -50908
$ws.Range("N129").Value = -166786672

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7814.881
$ws.Range("J132").Value = 9137.799999999999
$ws.Range("L132").Value = 27413.4
$ws.Range("N132").Value = -32473.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 26999.5
$ws.Range("H61").Value = 4426.1665
$ws.Range("I61").Value = 4639.25
$ws.Range("K61").Value = 4639.25
$ws.Range("M61").Value = -4437.25
$ws.Range("H100").Value = 5674
$ws.Range("I100").Value = 5451.1816
$ws.Range("K100").Value = 5451.1816
$ws.Range("M100").Value = -4910.1816
$ws.Range("H113").Value = 4426.1665
$ws.Range("I113").Value = 4639.25
$ws.Range("K113").Value = 4639.25
$ws.Range("M113").Value = -2469.25
$ws.Range("H136").Value = 6562.619
$ws.Range("I136").Value = 5093.4614
$ws.Range("K136").Value = 15280.3842
$ws.Range("M136").Value = -12730.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2000973
$ws.Range("I107").Value = 2857825.2
$ws.Range("K107").Value = 8573475.600000001
$ws.Range("M107").Value = -8571555.600000001
$ws.Range("H126").Value = 1982.7
$ws.Range("I126").Value = 1814.1111
$ws.Range("K126").Value = 5442.3333
$ws.Range("M126").Value = -2972.3333

